# Auto-generated: refresh market-data-derived profit columns (H-N) across all class sheets
# per scheduled-runner market data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1132.4
$ws.Range("I15").Value = 1132.4
$ws.Range("K15").Value = 3397.2
$ws.Range("M15").Value = -3228.2
$ws.Range("H98").Value = 2841.5
$ws.Range("I98").Value = 3009.8
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 3009.8
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -1511.8
$ws.Range("N98").Value = -4996
$ws.Range("H122").Value = 2841.5
$ws.Range("I122").Value = 3009.8
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9029.400000000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6579.400000000001
$ws.Range("N122").Value = -10900
$ws.Range("H137").Value = 1401.25
$ws.Range("I137").Value = 1607.4857
$ws.Range("J137").Value = 1290.2
$ws.Range("K137").Value = 4822.4571
$ws.Range("L137").Value = 3870.6
$ws.Range("M137").Value = -2272.4571
$ws.Range("N137").Value = -8970.6
$ws.Range("H138").Value = 5132146.5
$ws.Range("I138").Value = 2437
$ws.Range("J138").Value = 7412017.5
$ws.Range("K138").Value = 7311
$ws.Range("L138").Value = 22236052.5
$ws.Range("M138").Value = -2171
$ws.Range("N138").Value = -22246332.5
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 41745
$ws.Range("J140").Value = 41745
$ws.Range("L140").Value = 41745
$ws.Range("N140").Value = -52105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2903.203
$ws.Range("I32").Value = 2473.5334
$ws.Range("J32").Value = 5767.6665
$ws.Range("K32").Value = 2473.5334
$ws.Range("L32").Value = 5767.6665
$ws.Range("M32").Value = -2186.5334
$ws.Range("N32").Value = -6341.6665
$ws.Range("H63").Value = 3212.3076
$ws.Range("I63").Value = 2220
$ws.Range("J63").Value = 4800
$ws.Range("K63").Value = 2220
$ws.Range("L63").Value = 4800
$ws.Range("M63").Value = -1534
$ws.Range("N63").Value = -6172
$ws.Range("H66").Value = 3212.3076
$ws.Range("I66").Value = 2220
$ws.Range("J66").Value = 4800
$ws.Range("K66").Value = 11100
$ws.Range("L66").Value = 24000
$ws.Range("M66").Value = -7668
$ws.Range("N66").Value = -30864
$ws.Range("H74").Value = 27831.621
$ws.Range("I74").Value = 39036.152
$ws.Range("J74").Value = 1348.1818
$ws.Range("K74").Value = 39036.152
$ws.Range("L74").Value = 1348.1818
$ws.Range("M74").Value = -38162.152
$ws.Range("N74").Value = -3096.1818
$ws.Range("H77").Value = 27831.621
$ws.Range("I77").Value = 39036.152
$ws.Range("J77").Value = 1348.1818
$ws.Range("K77").Value = 195180.76
$ws.Range("L77").Value = 6740.909000000001
$ws.Range("M77").Value = -190812.76
$ws.Range("N77").Value = -15476.909
$ws.Range("H122").Value = 2100
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 531.6818
$ws.Range("I64").Value = 287.63635
$ws.Range("J64").Value = 775.7273
$ws.Range("K64").Value = 287.63635
$ws.Range("L64").Value = 775.7273
$ws.Range("M64").Value = -62.63634999999999
$ws.Range("N64").Value = -1225.7273
$ws.Range("H67").Value = 531.6818
$ws.Range("I67").Value = 287.63635
$ws.Range("J67").Value = 775.7273
$ws.Range("K67").Value = 287.63635
$ws.Range("L67").Value = 775.7273
$ws.Range("M67").Value = 492.36365
$ws.Range("N67").Value = -2335.7273
$ws.Range("H105").Value = 2239.7222
$ws.Range("I105").Value = 2105.95
$ws.Range("J105").Value = 2406.9375
$ws.Range("K105").Value = 2105.95
$ws.Range("L105").Value = 2406.9375
$ws.Range("M105").Value = -358.9499999999998
$ws.Range("N105").Value = -5900.9375
$ws.Range("H107").Value = 1979.9131
$ws.Range("I107").Value = 1979
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1979
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -59
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 11857
$ws.Range("I36").Value = 6500
$ws.Range("J36").Value = 13999.8
$ws.Range("K36").Value = 6500
$ws.Range("L36").Value = 13999.8
$ws.Range("M36").Value = -6112
$ws.Range("N36").Value = -14775.8
$ws.Range("H40").Value = 11857
$ws.Range("I40").Value = 6500
$ws.Range("J40").Value = 13999.8
$ws.Range("K40").Value = 6500
$ws.Range("L40").Value = 13999.8
$ws.Range("M40").Value = -6340
$ws.Range("N40").Value = -14319.8
$ws.Range("H107").Value = 1231.2307
$ws.Range("I107").Value = 595.8461
$ws.Range("J107").Value = 1866.6154
$ws.Range("K107").Value = 595.8461
$ws.Range("L107").Value = 1866.6154
$ws.Range("M107").Value = 1324.1539
$ws.Range("N107").Value = -5706.6154

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1506.95
$ws.Range("I68").Value = 802.0645
$ws.Range("J68").Value = 1952.898
$ws.Range("K68").Value = 2406.1935
$ws.Range("L68").Value = 5858.694
$ws.Range("M68").Value = -1595.1935
$ws.Range("N68").Value = -7480.694
$ws.Range("H71").Value = 1506.95
$ws.Range("I71").Value = 802.0645
$ws.Range("J71").Value = 1952.898
$ws.Range("K71").Value = 7218.5805
$ws.Range("L71").Value = 17576.082
$ws.Range("M71").Value = -3162.5805
$ws.Range("N71").Value = -25688.082
$ws.Range("H131").Value = 860.61
$ws.Range("I131").Value = 496.22223
$ws.Range("J131").Value = 896.6484
$ws.Range("K131").Value = 1488.66669
$ws.Range("L131").Value = 2689.9452
$ws.Range("M131").Value = 3551.33331
$ws.Range("N131").Value = -12769.9452
$ws.Range("H134").Value = 7588.852
$ws.Range("I134").Value = 4579.8
$ws.Range("J134").Value = 8272.727999999999
$ws.Range("K134").Value = 13739.4
$ws.Range("L134").Value = 24818.184
$ws.Range("M134").Value = -8669.400000000001
$ws.Range("N134").Value = -34958.18399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1467.1428
$ws.Range("I102").Value = 1498.4
$ws.Range("J102").Value = 1389
$ws.Range("K102").Value = 1498.4
$ws.Range("L102").Value = 1389
$ws.Range("M102").Value = 123.5999999999999
$ws.Range("N102").Value = -4633
$ws.Range("H132").Value = 2414.0557
$ws.Range("I132").Value = 2208.3914
$ws.Range("J132").Value = 2777.923
$ws.Range("K132").Value = 6625.174199999999
$ws.Range("L132").Value = 8333.769
$ws.Range("M132").Value = -4095.174199999999
$ws.Range("N132").Value = -13393.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3495
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 3990
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3990
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -4214
$ws.Range("H22").Value = 715.1667
$ws.Range("I22").Value = 225.75
$ws.Range("K22").Value = 225.75
$ws.Range("M22").Value = 69.25
$ws.Range("H27").Value = 715.1667
$ws.Range("I27").Value = 225.75
$ws.Range("K27").Value = 225.75
$ws.Range("M27").Value = -118.75
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8800
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3495
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3990
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 11970
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -16910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2214.75
$ws.Range("I81").Value = 1884.5
$ws.Range("J81").Value = 2324.8333
$ws.Range("K81").Value = 3769
$ws.Range("L81").Value = 4649.6666
$ws.Range("M81").Value = -2708
$ws.Range("N81").Value = -6771.6666
$ws.Range("H84").Value = 2214.75
$ws.Range("I84").Value = 1884.5
$ws.Range("J84").Value = 2324.8333
$ws.Range("K84").Value = 18845
$ws.Range("L84").Value = 23248.333
$ws.Range("M84").Value = -13541
$ws.Range("N84").Value = -33856.333
$ws.Range("H92").Value = 30525
$ws.Range("J92").Value = 30525
$ws.Range("L92").Value = 30525
$ws.Range("N92").Value = -35517
$ws.Range("H107").Value = 100002
$ws.Range("I107").Value = 100002
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 300006
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -298086
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 152224.38
$ws.Range("I122").Value = 2600
$ws.Range("J122").Value = 241999
$ws.Range("K122").Value = 7800
$ws.Range("L122").Value = 725997
$ws.Range("M122").Value = -5350
$ws.Range("N122").Value = -730897
$ws.Range("H126").Value = 71430190
$ws.Range("I126").Value = 1114.25
$ws.Range("J126").Value = 166668960
$ws.Range("K126").Value = 3342.75
$ws.Range("L126").Value = 500006880
$ws.Range("M126").Value = -872.75
$ws.Range("N126").Value = -500011820

